$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: add country ("Canada") in column I ---
$ws.Range("I7").Value = "Canada"

# --- Row 8: new occurrence record ---
$ws.Range("A8").Value = 37
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "Thysanocarpus curvipes"
$ws.Range("D8").Value = "Thysanocarpus curvipes"
$ws.Range("E8").Value = "h"
$ws.Range("F8").Value = "Thysanocarpus curvipes"
$ws.Range("G8").Value = 19750506

# Register "Mayne Island" (K8) before the longer locality string that also
# contains it (H8), so the new shared-string table order matches the
# original authoring order.
$ws.Range("K8").Value = "Mayne Island"
$ws.Range("H8").Value = "Mayne Island; cliffs facing navy channel"
$ws.Range("I8").Value = "Canada"
$ws.Range("J8").Value = "British Columbia"

# Row 8 grew tall because of the wrapped locality text - match the row height.
$ws.Rows.Item(8).RowHeight = 85

# --- Selection housekeeping to mirror the saved view state ---
$ws.Range("K8").Select()
